# Solved 2 questions from leetcode
#
# Adds a new tracker row (row 24) for the newly solved problem:
#   Category: "24. Sliding Window"
#   Name:     "Maximum Average Subarray I"
#   Link:     https://leetcode.com/problems/maximum-average-subarray-i/

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$category = "24. Sliding Window"
$name     = "Maximum Average Subarray I"
$url      = "https://leetcode.com/problems/maximum-average-subarray-i/"

# Duplicate the previous row's formatting (fill/font/hyperlink styles) into a
# freshly inserted row 24, the same way the rest of the sheet is laid out.
$ws.Rows(23).Copy()
$ws.Rows(24).Insert()

# Fill in the new row's values.
$ws.Cells.Item(24, 4).Value = $url
$ws.Cells.Item(24, 2).Value = $name
$ws.Cells.Item(24, 1).Value = $category

# Wire up the hyperlink for the Link column, matching the style used by the
# rest of the "Link" column.
$linkCell = $ws.Cells.Item(24, 4)
$ws.Hyperlinks.Add($linkCell, $url)
$linkCell.Style = "Hyperlink"

Write-Output "Added row 24 (Maximum Average Subarray I) to the tracker sheet"
